$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel
# (e.g. "1.00", "0.999", "16.00") are forced to Text via NumberFormat "@"
# before the value is written, then the style is reset back to Normal so no
# visible formatting/style change is introduced - only the text content changes.
$numericLookingCells = @("D4","D5","D11","D20","D21","D22","D24","D25","D27","D35","D41","D45","D46","D47","D50")
foreach ($ref in $numericLookingCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = '68.259.77'
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").Value = '3.563.30'
$ws.Range("E3").Value = '  +1.93%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '619.16'
$ws.Range("E5").Value = '  +3.05%  '
$ws.Range("E6").Value = '  +3.82%  '
$ws.Range("D7").Value = '3.562.44'
$ws.Range("E7").Value = '  +1.95%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +5.70%  '
$ws.Range("D11").Value = '7.47'
$ws.Range("E11").Value = '  +7.89%  '
$ws.Range("E12").Value = '  +3.97%  '
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D15").Value = '4.165.91'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '3.560.93'
$ws.Range("E16").Value = '  +1.80%  '
$ws.Range("D17").Value = '68.286.71'
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("E19").Value = '  +5.68%  '
$ws.Range("D20").Value = '16.00'
$ws.Range("E20").Value = '  +7.01%  '
$ws.Range("D21").Value = '10.04'
$ws.Range("E21").Value = '  +11.53%  '
$ws.Range("D22").Value = '454.16'
$ws.Range("E22").Value = '  +1.85%  '
$ws.Range("D24").Value = '78.47'
$ws.Range("D25").Value = '0.0000131'
$ws.Range("E25").Value = '  +2.69%  '
$ws.Range("D26").Value = '3.706.41'
$ws.Range("E26").Value = '  +1.92%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  +13.28%  '
$ws.Range("E29").Value = '  +4.15%  '
$ws.Range("E30").Value = '  +11.74%  '
$ws.Range("E31").Value = '  +3.69%  '
$ws.Range("E32").Value = '  +3.89%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +5.68%  '
$ws.Range("D35").Value = '26.13'
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("E36").Value = '  +5.04%  '
$ws.Range("D37").Value = '3.556.90'
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("E38").Value = '  +3.55%  '
$ws.Range("E39").Value = '  +8.81%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '181.66'
$ws.Range("E41").Value = '  +4.06%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  +4.92%  '
$ws.Range("D45").Value = '31.23'
$ws.Range("E45").Value = '  +14.87%  '
$ws.Range("D46").Value = '0.899'
$ws.Range("E46").Value = '  +2.29%  '
$ws.Range("D47").Value = '46.20'
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("E48").Value = '  +5.66%  '
$ws.Range("E49").Value = '  +4.45%  '
$ws.Range("D50").Value = '7.78'
$ws.Range("E50").Value = '  +3.56%  '
$ws.Range("E51").Value = '  +7.84%  '

# Restore default style on the cells we temporarily reformatted as Text
foreach ($ref in $numericLookingCells) {
    $ws.Range($ref).Style = "Normal"
}
